# Hjemme passive tweaks lichtwark deleted values
# Updates the "Lichtwark" sample-size / value columns (B:E) on rows 1-3 and
# narrows the active selection from B1:AY3 down to B1:E3 to reflect the
# columns that were actually re-computed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - sample counts for CON / STR measures (columns B-E)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 - "CON" data, recomputed after removing some passive values
$ws.Range("B2").Value = 8.2126824199167459
$ws.Range("C2").Value = 2.0103274370086979
$ws.Range("D2").Value = 1.1832059481919259
$ws.Range("E2").Value = 0.22662787183105593

# Row 3 - "STR" data, recomputed after removing some passive values
$ws.Range("B3").Value = 6.9053078809481292
$ws.Range("C3").Value = 12.606223091713815
$ws.Range("D3").Value = 9.02427344368121
$ws.Range("E3").Value = -4.2811632553267032

# Shrink the saved selection to match the updated range
$ws.Range("B1:E3").Select()
